$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -22.26
$ws.Range("A10").Value = -21.727
$ws.Range("A12").Value = -21.696
$ws.Range("A18").Value = -22.187
$ws.Range("A25").Value = -21.737
$ws.Range("A37").Value = -20.232
$ws.Range("A55").Value = -22.164
$ws.Range("A68").Value = -21.696
$ws.Range("A77").Value = -20.48
$ws.Range("A78").Value = -19.915
$ws.Range("A79").Value = -21.57
$ws.Range("A80").Value = -20.193
$ws.Range("A81").Value = -21.797
$ws.Range("A82").Value = -22.152
$ws.Range("A84").Value = -22.047
$ws.Range("A101").Value = -21.249
$ws.Range("A102").Value = -20.501
